$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Artn"
$ws.Range("C2").Value = "Gfra3"
$ws.Range("D2").Value = "Neutro"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.534048
$ws.Range("H2").Value = 1.602144
$ws.Range("I2").Value = 0.4424453651869046
$ws.Range("J2").Value = 0.5032206999227333
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.01376766666666667
$ws.Range("N2").Value = 0.041303
$ws.Range("O2").Value = 0.008447791226675409
$ws.Range("P2").Value = 0.01261838808592956
$ws.Range("Q2").Value = 0.007352594847999999
$ws.Range("R2").Value = 0.066173353632
$ws.Range("S2").Value = 0.003737686074309131
$ws.Range("T2").Value = 0.006349834084498155

$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Artn"
$ws.Range("C3").Value = "Gfra3"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.534048
$ws.Range("H3").Value = 1.602144
$ws.Range("I3").Value = 0.4424453651869046
$ws.Range("J3").Value = 0.5032206999227333
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.615968
$ws.Range("N3").Value = 3.231936
$ws.Range("O3").Value = 0.9915522087733246
$ws.Range("P3").Value = 0.9873816119140704
$ws.Range("Q3").Value = 0.863004478464
$ws.Range("R3").Value = 5.178026870784
$ws.Range("S3").Value = 0.4387076791125955
$ws.Range("T3").Value = 0.4968708658382351

$ws.Range("A4").Value = "M2"
$ws.Range("B4").Value = "Artn"
$ws.Range("C4").Value = "Gfra3"
$ws.Range("D4").Value = "Neutro"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.2356576666666667
$ws.Range("H4").Value = 0.706973
$ws.Range("I4").Value = 0.1952364626165198
$ws.Range("J4").Value = 0.2220546017626846
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.01376766666666667
$ws.Range("N4").Value = 0.041303
$ws.Range("O4").Value = 0.008447791226675409
$ws.Range("P4").Value = 0.01261838808592956
$ws.Range("Q4").Value = 0.003244456202111111
$ws.Range("R4").Value = 0.029200105819
$ws.Range("S4").Value = 0.001649316876018977
$ws.Range("T4").Value = 0.002801971141308093

$ws.Range("A5").Value = "M2"
$ws.Range("B5").Value = "Artn"
$ws.Range("C5").Value = "Gfra3"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.2356576666666667
$ws.Range("H5").Value = 0.706973
$ws.Range("I5").Value = 0.1952364626165198
$ws.Range("J5").Value = 0.2220546017626846
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.615968
$ws.Range("N5").Value = 3.231936
$ws.Range("O5").Value = 0.9915522087733246
$ws.Range("P5").Value = 0.9873816119140704
$ws.Range("Q5").Value = 0.380815248288
$ws.Range("R5").Value = 2.284891489728
$ws.Range("S5").Value = 0.1935871457405008
$ws.Range("T5").Value = 0.2192526306213765

$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Artn"
$ws.Range("C6").Value = "Gfra3"
$ws.Range("D6").Value = "Neutro"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.5
$ws.Range("G6").Value = 0.4373315
$ws.Range("H6").Value = 0.874663
$ws.Range("I6").Value = 0.3623181721965756
$ws.Range("J6").Value = 0.274724698314582
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.01376766666666667
$ws.Range("N6").Value = 0.041303
$ws.Range("O6").Value = 0.008447791226675409
$ws.Range("P6").Value = 0.01261838808592956
$ws.Range("Q6").Value = 0.006021034314833333
$ws.Range("R6").Value = 0.036126205889
$ws.Range("S6").Value = 0.003060788276347301
$ws.Range("T6").Value = 0.003466582860123315

$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Artn"
$ws.Range("C7").Value = "Gfra3"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.5
$ws.Range("G7").Value = 0.4373315
$ws.Range("H7").Value = 0.874663
$ws.Range("I7").Value = 0.3623181721965756
$ws.Range("J7").Value = 0.274724698314582
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.615968
$ws.Range("N7").Value = 3.231936
$ws.Range("O7").Value = 0.9915522087733246
$ws.Range("P7").Value = 0.9873816119140704
$ws.Range("Q7").Value = 0.706713709392
$ws.Range("R7").Value = 2.826854837568
$ws.Range("S7").Value = 0.3592573839202283
$ws.Range("T7").Value = 0.2712581154544587
